# 20171221 修正Google Excel問題單 第78 80 99條
# Add a new "改單人員" (U_NAME) column (J) to the 組長歷史查詢結果 report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell J1 ("改單人員") ---
# Copy formatting from I1 (the last existing header cell) so the new header
# cell gets the same style (bold header style), then set its text.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "改單人員"

# --- New data cell J2 ("${table:data.U_NAME}") ---
# Copy formatting from I2 (the last existing data cell) so the new data
# cell gets the same style, then set its placeholder text.
$ws.Range("I2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$ws.Range("J2").Value = '${table:data.U_NAME}'

$excel.CutCopyMode = $false

# Match the column width of the neighbouring data columns (E:I = 13.5 chars).
$ws.Columns.Item(10).ColumnWidth = 13.5

# Update the active selection as it was left after the edit.
[void]$ws.Range("F3").Select()
